$d = $word.ActiveDocument

# Locate the "Add Blog.find()" paragraph which is the anchor point for the new content
# (new SHOW ROUTE / EDIT ROUTE sections are inserted right after it, before the blank
# paragraph that precedes "THINGS TO RESEARCH").
$anchorRange = $d.Content
$found = $anchorRange.Find.Execute( `
    "Add Blog.find() so the current blogs from the database are returned and can be accessed", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph"
}
$anchorPara = $anchorRange.Paragraphs(1)
$anchorIndex = $anchorPara.Index

# --- empty paragraph ---
$anchorPara.Range.InsertParagraphAfter()
$p = $d.Paragraphs($anchorIndex + 1)
$p.Style = "Normal"

# --- "SHOW ROUTE" heading ---
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($anchorIndex + 2)
$p.Style = "Heading 5"
$p.Format.SpaceBefore = 0
$p.Range.InsertAfter("SHOW ROUTE")

# --- "- Create show.ejs file" ---
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($anchorIndex + 3)
$p.Style = "Normal"
$p.Range.InsertAfter("- Create show.ejs file")

# --- "- Create show route which will access the specific show page based off its _id from the database" ---
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($anchorIndex + 4)
$p.Style = "Normal"
$p.Range.InsertAfter("- Create show route which will access the specific show page based off its _id from the database")

# --- empty paragraph ---
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($anchorIndex + 5)
$p.Style = "Normal"

# --- "EDIT ROUTE" heading; the _GoBack bookmark moves here, between "EDIT" and " ROUTE" ---
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($anchorIndex + 6)
$p.Style = "Heading 5"
$p.Format.SpaceBefore = 0
$editStart = $p.Range.Start
$p.Range.InsertAfter("EDIT ROUTE")
$bmPos = $editStart + 4
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

# --- "- Added update route for blog" ---
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($anchorIndex + 7)
$p.Style = "Normal"
$p.Range.InsertAfter("- Added update route for blog")

# --- "- Added edit page with form" ---
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($anchorIndex + 8)
$p.Style = "Normal"
$p.Range.InsertAfter("- Added edit page with form")

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
